# Applies the diff described in the commit:
# "Added a line of code to convert all values to float after being added to
#  dataframe. This will ensure all figures have 2 decimal places when later
#  code is executed" — plus a new "HD" position for Michael B that flows
# through the Transactions and Summary sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Investors"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Investors")
$ws.Range("A2").Value = "Summary"
$ws.Range("B2:B5").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Sheet "Summary"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Summary")

# Row 2 - NVDA
$ws.Range("B2").Value = 809.75
$ws.Range("E2").Value = 11336.5
$ws.Range("H2").Value = 7304.5

# Row 3 - MSTR
$ws.Range("B3").Value = 1196.54
$ws.Range("D3").Value = 0
$ws.Range("D3").NumberFormat = "0.00"
$ws.Range("E3").Value = 5982.7
$ws.Range("H3").Value = 3588.7

# Row 4 - COIN
$ws.Range("B4").Value = 215.52
$ws.Range("E4").Value = 646.5599999999999
$ws.Range("H4").Value = 468.3

# Row 5 - CCOR
$ws.Range("B5").Value = 26.44
$ws.Range("D5").Value = 0
$ws.Range("D5").NumberFormat = "0.00"
$ws.Range("E5").Value = 132.2
$ws.Range("H5").Value = 24.6

# Row 6 - INTC
$ws.Range("B6").Value = 34.53
$ws.Range("D6").Value = 0
$ws.Range("D6").NumberFormat = "0.00"
$ws.Range("E6").Value = 517.95
$ws.Range("H6").Value = -281.7

# Row 7 - MTB
$ws.Range("B7").Value = 142.45
$ws.Range("D7").Value = 0
$ws.Range("D7").NumberFormat = "0.00"
$ws.Range("E7").Value = 1709.4
$ws.Range("H7").Value = 28.8

# Row 8 - AMD
$ws.Range("B8").Value = 148.07
$ws.Range("D8").Value = 0
$ws.Range("D8").NumberFormat = "0.00"
$ws.Range("E8").Value = 740.35
$ws.Range("H8").Value = -34.65

# Row 9 - NEW - HD
$ws.Range("A9").Value = "HD"
$ws.Range("B9").Value = 334.53
$ws.Range("C9").Value = 12
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 4014.36
$ws.Range("F9").Value = 4000.68
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 13.68
$ws.Range("B9:H9").NumberFormat = "0.00"

# Make sure all the pre-existing numeric cells on rows 2-8 that didn't change
# value still pick up the new "convert to float" style.
$ws.Range("B2:H8").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Sheet "Transactions"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Transactions")

$ws.Range("B2:B10").NumberFormat = "0.00"
$ws.Range("F2:F10").NumberFormat = "0.00"
$ws.Range("H3").Value = 799.65

# Row 11 - NEW transaction - Michael B buys HD
$ws.Range("A11").Value = "Michael B"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 45401
$ws.Range("C11").NumberFormat = "YYYY-MM-DD"
$ws.Range("D11").Value = "HD"
$ws.Range("E11").Value = "Buy"
$ws.Range("F11").Value = 12
$ws.Range("G11").Value = 333.39
$ws.Range("H11").Value = 4000.68
$ws.Range("B11").NumberFormat = "0.00"
$ws.Range("F11:H11").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Sheet "Joe L"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Joe L")
$ws.Range("C2:C4").NumberFormat = "0.00"
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D3:D4").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Sheet "Jonathan R"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Jonathan R")
$ws.Range("C2:C5").NumberFormat = "0.00"
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D4:D5").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Sheet "Michael B"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Michael B")
$ws.Range("D2").Value = 0
$ws.Range("B2:H2").NumberFormat = "0.00"

Write-Output "edit applied"
